# Update the "variable" column (E) labels to remove the "Xnn.01kg.a.nnkg"
# prefix text, keeping only the trailing weight-bracket number, per the
# freight table (carga fracionada) fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E14").Value = "30"
$ws.Range("E15:E27").Value = "50"
$ws.Range("E28:E40").Value = "70"
$ws.Range("E41:E53").Value = "100"
